$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.762.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.094.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.37%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'228.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'60.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.88%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.73%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.55%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.43%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'15.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.98%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.404.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'21.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.30%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.819"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +6.48%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.065.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'38.687.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'71.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0841"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'227.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.43%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.54%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Cosmos"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'9.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.24%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Monero"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'171.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.21%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +7.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +10.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'19.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.52%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'WEMIXToken"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'2.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.69%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Stellar"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.56%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.99%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0231"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.76%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Aave"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'101.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Maker"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.537.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.66%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0919"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +8.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.30%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.63%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.290.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.32%  "
$ws.Range("E51").Style = "Normal"
